$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version bump: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date refresh
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value now populated
$ws.Range("B9").Value = "Alvearie Team"

# Remove the duplicated "Contact" / "No display for ContactDetail" row
# (row 11 was an exact duplicate of row 10); deleting it shifts rows 12-22 up.
$ws.Rows.Item(11).Delete()

# Row 10 (formerly "Contact") becomes "Jurisdiction"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# "Case Sensitive" value now populated (row 14 after the deletion above).
# Build it via a formula then collapse to a literal value so it lands as
# text "true" rather than being auto-coerced into a boolean.
$ws.Range("B14").Formula = "=""tru""&""e"""
$ws.Range("B14").Copy()
$ws.Range("B14").PasteSpecial(-4163)
